$d = $word.ActiveDocument

# Locate the paragraph that ends with "/client" (the "remixd -s c:/bccode/retoken-app/client" line)
$rng = $d.Content
$rng.Find.Execute("/client", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para = $rng.Paragraphs(1)

# Insert a brand-new paragraph right after it (before the existing blank paragraph)
$prange = $para.Range
$prange.Collapse(0)
$prange.InsertParagraphAfter()
$newPara = $para.Next()
$nrange = $newPara.Range

# Build the new paragraph's OOXML exactly (Courier New run formatting + proofErr
# spell-check markers around "remixd", matching the q-text command-line style)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="q-text"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="282829"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="282829"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">$ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="282829"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>remixd</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="282829"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> -s c:/bccode/retoken-app</w:t></w:r></w:p>'

$nrange.InsertXML($xml)

$d.Save()
